# Update the division-problem table cells in place.
# The table has data only in rows 1, 5, 9, 13, 17 (1-indexed),
# with 5 columns of "N÷N=" expressions each.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Row => (col => new text) mapping for each data-bearing table row.
$updates = @{
    1  = @("80÷4=", "50÷2=", "88÷5=", "79÷9=", "79÷3=")
    5  = @("87÷5=", "72÷8=", "55÷3=", "96÷4=", "58÷2=")
    9  = @("24÷2=", "90÷8=", "41÷4=", "58÷6=", "16÷2=")
    13 = @("39÷6=", "86÷7=", "45÷8=", "87÷7=", "23÷5=")
    17 = @("32÷3=", "26÷4=", "22÷3=", "13÷3=", "47÷4=")
}

foreach ($rowIndex in $updates.Keys) {
    $newValues = $updates[$rowIndex]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $table.Cell($rowIndex, $col)
        $range = $cell.Range
        # Trim the trailing cell-mark/paragraph-mark characters Word appends
        # to cell ranges so only the visible text is replaced.
        $range.End = $range.End - 1
        $range.Text = $newValues[$col - 1]
    }
}
